$d = $word.ActiveDocument

# Part 2.1: remove the "10%" tax rate mention and replace it with "в США",
# and drop the now-redundant "(3%) " qualifier on "НДФЛ".
$d.Content.Find.Execute(
    "удерживался налог по ставке 10%, поэтому",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "удерживался налог в США, поэтому", 2)

$d.Content.Find.Execute(
    "неоплаченной суммы НДФЛ (3%) в рублях",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "неоплаченной суммы НДФЛ в рублях", 2)
